$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (Changed) date column C for all existing data rows
#    (rows 2..526) from 2023-09-13 (45182) to 2023-09-15 (45184).
$ws.Range("C2:C526").Value = 45184

# 2. Make row 526's height explicit (ht="15" customHeight="1"), matching the
#    rest of the sheet's rows.
$ws.Rows.Item(526).RowHeight = 15

# 3. Append new row 527.
$ws.Range("A527").Value = "A 42899-2023"
$ws.Range("B527").Value = 45182
$ws.Range("B527").NumberFormat = "YYYY-MM-DD"
$ws.Range("C527").Value = 45184
$ws.Range("C527").NumberFormat = "YYYY-MM-DD"
$ws.Range("D527").Value = "SKÅNE LÄN"
$ws.Range("E527").Value = "KRISTIANSTAD"
$ws.Range("F527").Value = "Sveaskog"
$ws.Range("G527").Value = 1.5
$ws.Range("H527").Value = 0
$ws.Range("I527").Value = 0
$ws.Range("J527").Value = 0
$ws.Range("K527").Value = 0
$ws.Range("L527").Value = 0
$ws.Range("M527").Value = 0
$ws.Range("N527").Value = 0
$ws.Range("O527").Value = 0
$ws.Range("P527").Value = 0
$ws.Range("Q527").Value = 0
$ws.Range("R527").WrapText = $true
$ws.Rows.Item(527).RowHeight = 15

# 4. Append new row 528 (note: no Markägare/F value, and row height left at
#    the sheet default, i.e. no explicit customHeight attribute).
$ws.Range("A528").Value = "A 43321-2023"
$ws.Range("B528").Value = 45183
$ws.Range("B528").NumberFormat = "YYYY-MM-DD"
$ws.Range("C528").Value = 45184
$ws.Range("C528").NumberFormat = "YYYY-MM-DD"
$ws.Range("D528").Value = "SKÅNE LÄN"
$ws.Range("E528").Value = "KRISTIANSTAD"
$ws.Range("G528").Value = 1
$ws.Range("H528").Value = 0
$ws.Range("I528").Value = 0
$ws.Range("J528").Value = 0
$ws.Range("K528").Value = 0
$ws.Range("L528").Value = 0
$ws.Range("M528").Value = 0
$ws.Range("N528").Value = 0
$ws.Range("O528").Value = 0
$ws.Range("P528").Value = 0
$ws.Range("Q528").Value = 0
$ws.Range("R528").WrapText = $true
